# Log_of_all_Blogs.xlsx - add Post 72 entry to the blog log table (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Expand the existing table ("Table2") by one row; this grows the table
# range / autofilter range and creates a new underlying worksheet row.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# Carry over the formatting (date format, hyperlink look, borders, etc.)
# from the previous row so the new row matches the rest of the table.
$ws.Range("B81:F81").Copy()
$ws.Range("B82:F82").PasteSpecial(-4122)

# Fill in the new row's cells. Order matters for the resulting shared
# string table layout (dev.to link, then title, then hashnode link).
$ws.Range("B82").Value2 = 72
$ws.Range("F82").Value2 = "https://dev.to/rahulmishra05/various-allocation-methods-in-contiguous-memory-management-operating-system-m05-p05-1k1c"
$ws.Range("C82").Value2 = "Various Allocation Methods in Contiguous Memory Management | Operating System - M05 P05"
$ws.Range("D82").Value2 = 44184
$ws.Range("E82").Value2 = "https://programmingport.hashnode.dev/various-allocation-methods-in-contiguous-memory-management-or-operating-system-m05-p05"

# Widen column C (Title of the Blog Post) to fit the new, longer title.
$ws.Columns.Item(3).ColumnWidth = 80.6666666666667

# Update the view: scroll so the new row is visible, and select the new
# hashnode-link cell, mirroring where the author's cursor ended up.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E82").Select() | Out-Null
